# Auto-generated Excel COM-interop script
# Applies the numeric cell updates (set / add / remove) described by the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2783714.5
$ws.Range("M116").Value = -2192.857
$ws.Range("K116").Value = 5634.857
$ws.Range("I116").Value = 5634.857
$ws.Range("I132").Value = 1490.421
$ws.Range("K132").Value = 4471.263
$ws.Range("M132").Value = -1941.263
$ws.Range("H132").Value = 1889.7805
$ws.Range("H138").Value = 1551.4286
$ws.Range("N138").Value = -19277
$ws.Range("M138").Value = 1670.0908
$ws.Range("I138").Value = 1156.6364
$ws.Range("J138").Value = 2999
$ws.Range("K138").Value = 3469.9092
$ws.Range("L138").Value = 8997
$ws.Range("M141").Value = -3557.6158
$ws.Range("H141").Value = 3390.8667
$ws.Range("I141").Value = 2912.5386
$ws.Range("K141").Value = 8737.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N2").Value = -1226
$ws.Range("J2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("M44").Value = -48178.668
$ws.Range("K44").Value = 48666.668
$ws.Range("J44").Value = 86144
$ws.Range("N44").Value = -87120
$ws.Range("H44").Value = 63657.6
$ws.Range("L44").Value = 86144
$ws.Range("I44").Value = 48666.668
$ws.Range("J62").Value = 0
$ws.Range("H62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("L64").Value = 0
$ws.Range("K64").Value = 30000
$ws.Range("H64").Value = 30000
$ws.Range("M64").Value = -29752
$ws.Range("N64").Value = ""
$ws.Range("J64").Value = 0
$ws.Range("I64").Value = 30000
$ws.Range("J65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("I67").Value = 30000
$ws.Range("M67").Value = -29142
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 30000
$ws.Range("H67").Value = 30000
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""
$ws.Range("L97").Value = 3000
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 510.70587
$ws.Range("I97").Value = 510.70587
$ws.Range("N97").Value = -3992
$ws.Range("H97").Value = 649
$ws.Range("M97").Value = -14.70587
$ws.Range("K102").Value = 49321.24
$ws.Range("I102").Value = 49321.24
$ws.Range("M102").Value = -47699.24
$ws.Range("H102").Value = 48178.457
$ws.Range("N110").Value = -5788
$ws.Range("H110").Value = 1799.3334
$ws.Range("L110").Value = 1698
$ws.Range("J110").Value = 1698
$ws.Range("L116").Value = 1000
$ws.Range("N116").Value = -5588
$ws.Range("J116").Value = 1000
$ws.Range("L128").Value = 69429
$ws.Range("J128").Value = 69429
$ws.Range("N128").Value = -79389
$ws.Range("H128").Value = 69429
$ws.Range("H138").Value = 40428.5
$ws.Range("N138").Value = -50708.5
$ws.Range("J138").Value = 40428.5
$ws.Range("L138").Value = 40428.5
$ws.Range("H140").Value = 66327.664
$ws.Range("N140").Value = -76687.664
$ws.Range("L140").Value = 66327.664
$ws.Range("J140").Value = 66327.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N3").Value = -1228
$ws.Range("J3").Value = 1000
$ws.Range("L3").Value = 1000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -5927
$ws.Range("H86").Value = 7050
$ws.Range("J86").Value = 0
$ws.Range("N86").Value = ""
$ws.Range("K86").Value = 7050
$ws.Range("I86").Value = 7050
$ws.Range("M89").Value = -29634
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("K89").Value = 35250
$ws.Range("I89").Value = 7050
$ws.Range("H89").Value = 7050
$ws.Range("N89").Value = ""
$ws.Range("I94").Value = 1613.5161
$ws.Range("K94").Value = 1613.5161
$ws.Range("M94").Value = -1162.5161
$ws.Range("H94").Value = 1666.9459
$ws.Range("I99").Value = 144544
$ws.Range("H99").Value = 2689567.2
$ws.Range("K99").Value = 144544
$ws.Range("M99").Value = -143046
$ws.Range("N126").Value = -74870
$ws.Range("J126").Value = 64990
$ws.Range("H126").Value = 64990
$ws.Range("L126").Value = 64990

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").Value = -2378.8
$ws.Range("H31").Value = 4037.4707
$ws.Range("K31").Value = 2673.8
$ws.Range("I31").Value = 2673.8
$ws.Range("M34").Value = -2471.8
$ws.Range("I34").Value = 2673.8
$ws.Range("H34").Value = 4037.4707
$ws.Range("K34").Value = 2673.8
$ws.Range("M44").Value = -1558
$ws.Range("K44").Value = 2000
$ws.Range("H44").Value = 2000
$ws.Range("I44").Value = 2000
$ws.Range("H58").Value = 1501.0416
$ws.Range("M58").Value = -1229.579
$ws.Range("K58").Value = 1432.579
$ws.Range("I58").Value = 1432.579
$ws.Range("K62").Value = 3900
$ws.Range("I62").Value = 3900
$ws.Range("J62").Value = 3350.5
$ws.Range("H62").Value = 3625.25
$ws.Range("L62").Value = 3350.5
$ws.Range("N62").Value = -4598.5
$ws.Range("M62").Value = -3276
$ws.Range("L63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("N63").Value = ""
$ws.Range("J63").Value = 0
$ws.Range("J65").Value = 3350.5
$ws.Range("M65").Value = -16380
$ws.Range("I65").Value = 3900
$ws.Range("K65").Value = 19500
$ws.Range("H65").Value = 3625.25
$ws.Range("L65").Value = 16752.5
$ws.Range("N65").Value = -22992.5
$ws.Range("L66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("N66").Value = ""
$ws.Range("H66").Value = 0
$ws.Range("J99").Value = 7814669
$ws.Range("I99").Value = 11113213
$ws.Range("H99").Value = 10170772
$ws.Range("L99").Value = 7814669
$ws.Range("K99").Value = 11113213
$ws.Range("N99").Value = -7817665
$ws.Range("M99").Value = -11111715
$ws.Range("N126").Value = -23448947
$ws.Range("I126").Value = 11113213
$ws.Range("J126").Value = 7814669
$ws.Range("H126").Value = 10170772
$ws.Range("K126").Value = 33339639
$ws.Range("M126").Value = -33337169
$ws.Range("L126").Value = 23444007
$ws.Range("I134").Value = 3108296.5
$ws.Range("H134").Value = 2500051.2
$ws.Range("K134").Value = 9324889.5
$ws.Range("M134").Value = -9322354.5
$ws.Range("H136").Value = 1501.0416
$ws.Range("M136").Value = -1747.737
$ws.Range("K136").Value = 4297.737
$ws.Range("I136").Value = 1432.579

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L64").Value = 87799
$ws.Range("H64").Value = 81399.5
$ws.Range("N64").Value = -88295
$ws.Range("J64").Value = 87799
$ws.Range("J67").Value = 87799
$ws.Range("H67").Value = 81399.5
$ws.Range("L67").Value = 87799
$ws.Range("N67").Value = -89515
$ws.Range("K97").Value = 360
$ws.Range("I97").Value = 360
$ws.Range("H97").Value = 466.66666
$ws.Range("M97").Value = 136
$ws.Range("I132").Value = 6548.6665
$ws.Range("K132").Value = 19645.9995
$ws.Range("M132").Value = -17115.9995
$ws.Range("H132").Value = 6981.091

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2201.5557
$ws.Range("I46").Value = 1577.6666
$ws.Range("M46").Value = -1389.6666
$ws.Range("K46").Value = 1577.6666
$ws.Range("K122").Value = 428584260
$ws.Range("H122").Value = 100004550
$ws.Range("M122").Value = -428581810
$ws.Range("I122").Value = 142861420
$ws.Range("J132").Value = 0
$ws.Range("N132").Value = ""
$ws.Range("L132").Value = 0
$ws.Range("H132").Value = 3260.1765
$ws.Range("H136").Value = 1751.9667
$ws.Range("M136").Value = -1894.0434
$ws.Range("J136").Value = 2641.1428
$ws.Range("K136").Value = 4444.0434
$ws.Range("I136").Value = 1481.3478
$ws.Range("N136").Value = -13023.4284
$ws.Range("L136").Value = 7923.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K62").Value = 4142.4287
$ws.Range("I62").Value = 4142.4287
$ws.Range("H62").Value = 4142.4287
$ws.Range("M62").Value = -3518.4287
$ws.Range("M65").Value = -17592.1435
$ws.Range("I65").Value = 4142.4287
$ws.Range("K65").Value = 20712.1435
$ws.Range("H65").Value = 4142.4287
$ws.Range("L86").Value = 48999
$ws.Range("H86").Value = 48999
$ws.Range("J86").Value = 48999
$ws.Range("N86").Value = -51245
$ws.Range("J89").Value = 48999
$ws.Range("L89").Value = 244995
$ws.Range("H89").Value = 48999
$ws.Range("N89").Value = -256227
$ws.Range("I107").Value = 1290.3334
$ws.Range("K107").Value = 3871.0002
$ws.Range("M107").Value = -1951.0002
$ws.Range("H107").Value = 1884.1818
$ws.Range("H136").Value = 1973.9
$ws.Range("M136").Value = -1859.571599999999
$ws.Range("K136").Value = 4409.571599999999
$ws.Range("I136").Value = 1469.8572

